$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("forms")

# Update the send-email element-definitions JSON (row 3, column B) first so the
# shared-string table regenerates with "send-email" ident ordered before "address-form":
# dataType field removed from every element, and ELEM-54758 renamed to templParams.
$ws.Range("B3").Value = '{"ident":"send-email","elems":[{"ident":"ELEM-84079","elementType":"TITLE","texts":{"label":{"textType":"LABEL","texts":{"EN":"Send Email Task","DE":"Send Email Task"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"INPUT_TYPE","value":"text"},{"extraProp":"SIZE","value":"20"},{"extraProp":"SIZE_CLASS","value":"huge"}]},"value":"","required":false,"inline":false,"readOnly":false,"layoutWide":"SIXTEEN","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"smtpConfig","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"SMTP Configuration","DE":"SMTP Konfiguration"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"INPUT_TYPE","value":"text"},{"extraProp":"SIZE","value":"20"},{"extraProp":"SIZE_CLASS","value":"big"}]},"value":"default","required":false,"inline":false,"readOnly":false,"layoutWide":"EIGHT","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"subject","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"Subject","DE":"Betreff"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"INPUT_TYPE","value":"text"},{"extraProp":"SIZE","value":"20"}]},"value":"I18N Email Test-Subject: {{mail.info.subject}}!","required":true,"inline":false,"readOnly":false,"layoutWide":"SIXTEEN","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"content","elementType":"TEXTAREA","texts":{"label":{"textType":"LABEL","texts":{"EN":"Content","DE":"Inhalt"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"INPUT_TYPE","value":"text"},{"extraProp":"SIZE","value":"20"}]},"value":"Guten Tag\n\nDies ist ein Test Email.\n\nFreundliche Grüsse\n\nBPF-Team","required":false,"inline":false,"readOnly":false,"layoutWide":"SIXTEEN","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"from","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"From (Sender)","DE":"Von (Sender)"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"Email Address of the Sender.","DE":"Email Adresse des Senders."}}},"extras":{"propValues":[{"extraProp":"INPUT_TYPE","value":"text"},{"extraProp":"SIZE","value":"20"}]},"value":"sender@finnova.ch","required":true,"inline":false,"readOnly":false,"layoutWide":"EIGHT","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":true,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"to","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"To (Recipients)","DE":"To (Empfänger)"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"Semicolon separated list of all recipient''s email addresses .","DE":"Strichpunkt getrennte Liste aller Email-Adressen der Empfänger."}}},"extras":{"propValues":[{"extraProp":"INPUT_TYPE","value":"text"},{"extraProp":"SIZE","value":"20"}]},"value":"hans.muster@test.com; vreny.mueller@test.com","required":true,"inline":false,"readOnly":false,"layoutWide":"EIGHT","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"ELEM-73603","elementType":"DIVIDER","texts":{"label":{"textType":"LABEL","texts":{"EN":"Additional Attributes for HTML Emails","DE":"Weitere Werte für die HTML Emails"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"INPUT_TYPE","value":"text"},{"extraProp":"SIZE","value":"20"},{"extraProp":"SIZE_CLASS","value":"big"}]},"value":"","required":false,"inline":false,"readOnly":false,"layoutWide":"SIXTEEN","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"templName","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"Template Name","DE":"Template Name"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"The name of the HTML Template File.","DE":"Der Name der HTML Template Date."}}},"extras":{"propValues":[{"extraProp":"INPUT_TYPE","value":"text"},{"extraProp":"SIZE","value":"20"}]},"value":"info_mail.html","required":false,"inline":false,"readOnly":false,"layoutWide":"TEN","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"language","elementType":"DROPDOWN","texts":{"label":{"textType":"LABEL","texts":{"EN":"Language","DE":"Sprache"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"INPUT_TYPE","value":"text"},{"extraProp":"SIZE","value":"20"},{"extraProp":"CLEARABLE","value":"false"}]},"value":"DE","required":false,"inline":false,"readOnly":false,"layoutWide":"SIX","elemEntries":{"entries":[{"key":"DE","label":{"textType":"LABEL","texts":{"EN":"German","DE":"Deutsch"}},"ident":"ENTRY-3318"},{"key":"EN","label":{"textType":"LABEL","texts":{"EN":"English","DE":"Englisch"}},"ident":"ENTRY-3546"}]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"templParams","elementType":"TEXTAREA","texts":{"label":{"textType":"LABEL","texts":{"EN":"Template Parameters","DE":"Template Parameter"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"Parameters as a Map that are used in the Template.","DE":"Ein JSON-Objekt, mit den Parameter, welche dem Template mitgegeben werden."}}},"extras":{"propValues":[{"extraProp":"INPUT_TYPE","value":"text"},{"extraProp":"SIZE","value":"20"}]},"value":"{\n \"header\": \"Finnova rocks\",\n \"greeting\": \"Dear Mister President\",\n \"content\": \"Lory Logsi Dogsi Email content with Link: link:nextTask\",\n \"footer\": \"visit us!\",\n \"nextTask\": \"Go to next Task.\",\n \"nextTaskLink\": \"link:nextTask\",\n \"startInstance\": \"Start new Process Instance\",\n \"startInstanceLink\": \"link:startInstance\"\n }","required":false,"inline":false,"readOnly":false,"layoutWide":"SIXTEEN","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}}]}'

# Update the address-form element-definitions JSON (row 2, column B):
# dataType field removed from every element, and the TEXTFIELD-* element
# idents are renamed to match the data/mapping idents (street, number, postcode, town).
$ws.Range("B2").Value = '{"ident":"address-form","elems":[{"ident":"TITLE-60664","elementType":"TITLE","texts":{"label":{"textType":"LABEL","texts":{"EN":"Address","DE":"Adresse"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"SIZE_CLASS","value":"huge"}]},"value":"","required":false,"inline":false,"readOnly":false,"layoutWide":"SIXTEEN","elemEntries":{"entries":[]},"validations":{"rules":[]}},{"ident":"street","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"Street","DE":"Strasse"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"SIZE","value":"20"},{"extraProp":"INPUT_TYPE","value":"text"}]},"value":"Sonnenweg","required":true,"inline":false,"readOnly":false,"layoutWide":"TWELVE","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"number","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"Number","DE":"Nummer"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"SIZE","value":"20"},{"extraProp":"INPUT_TYPE","value":"text"}]},"value":"23a","required":true,"inline":false,"readOnly":false,"layoutWide":"FOUR","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"postcode","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"Postcode","DE":"Plz"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"SIZE","value":"20"},{"extraProp":"INPUT_TYPE","value":"number"}]},"value":"6414","required":true,"inline":false,"readOnly":false,"layoutWide":"FOUR","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":true,"params":{"intParam1":1000,"intParam2":9999}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"town","elementType":"TEXTFIELD","texts":{"label":{"textType":"LABEL","texts":{"EN":"Town","DE":"Ort"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"SIZE","value":"20"},{"extraProp":"INPUT_TYPE","value":"text"}]},"value":"Oberarth","required":true,"inline":false,"readOnly":false,"layoutWide":"TWELVE","elemEntries":{"entries":[]},"validations":{"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}}]}'

# Narrow column B (it used to be sized for the long JSON blobs).
$ws.Columns.Item(2).ColumnWidth = 32.666666666666664

# Move the saved selection/active cell to B11.
$ws.Activate() | Out-Null
$ws.Range("B11").Select() | Out-Null

